$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 58,5
$data[0,0] = '长沙特来电飞狐四方坪西区充电站'
$data[0,1] = '701号直流'
$data[0,2] = 45927.457337962966
$data[0,3] = 45936.347141203703
$data[0,4] = 213.35527777770767
$data[1,0] = '长沙特来电飞狐四方坪西区充电站'
$data[1,1] = '502号直流'
$data[1,2] = 45930.238043981481
$data[1,3] = 45936.347141203703
$data[1,4] = 146.61833333334653
$data[2,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[2,1] = '112号直流'
$data[2,2] = 45930.517060185186
$data[2,3] = 45936.347141203703
$data[2,4] = 139.92194444441702
$data[3,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[3,1] = '111号直流'
$data[3,2] = 45930.618518518517
$data[3,3] = 45936.347141203703
$data[3,4] = 137.48694444447756
$data[4,0] = '长沙特来电飞狐四方坪东区充电站'
$data[4,1] = '201号直流'
$data[4,2] = 45931.575543981482
$data[4,3] = 45936.347141203703
$data[4,4] = 114.5183333333116
$data[5,0] = '长沙特来电飞狐四方坪西区充电站'
$data[5,1] = '603号直流'
$data[5,2] = 45932.081099537034
$data[5,3] = 45936.347141203703
$data[5,4] = 102.38500000006752
$data[6,0] = '长沙特来电飞狐四方坪南区充电站'
$data[6,1] = '406号直流'
$data[6,2] = 45933.039143518516
$data[6,3] = 45936.347141203703
$data[6,4] = 79.391944444505498
$data[7,0] = '长沙特来电飞狐四方坪南区充电站'
$data[7,1] = '103号直流'
$data[7,2] = 45933.305023148147
$data[7,3] = 45936.347141203703
$data[7,4] = 73.010833333362825
$data[8,0] = '长沙特来电飞狐四方坪西区充电站'
$data[8,1] = '503号直流'
$data[8,2] = 45934.068425925929
$data[8,3] = 45936.347141203703
$data[8,4] = 54.689166666590609
$data[9,0] = '长沙特来电飞狐四方坪西区充电站'
$data[9,1] = '505号直流'
$data[9,2] = 45934.284421296295
$data[9,3] = 45936.347141203703
$data[9,4] = 49.505277777789161
$data[10,0] = '长沙特来电飞狐四方坪东区充电站'
$data[10,1] = '903号直流'
$data[10,2] = 45934.542071759257
$data[10,3] = 45936.347141203703
$data[10,4] = 43.321666666714009
$data[11,0] = '长沙特来电飞狐四方坪西区充电站'
$data[11,1] = '801号直流'
$data[11,2] = 45934.554988425924
$data[11,3] = 45936.347141203703
$data[11,4] = 43.011666666716337
$data[12,0] = '长沙特来电飞狐四方坪东区充电站'
$data[12,1] = '001B号直流'
$data[12,2] = 45934.55914351852
$data[12,3] = 45936.347141203703
$data[12,4] = 42.911944444407709
$data[13,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[13,1] = '110号直流'
$data[13,2] = 45934.707303240742
$data[13,3] = 45936.347141203703
$data[13,4] = 39.356111111061182
$data[14,0] = '长沙特来电飞狐四方坪西区充电站'
$data[14,1] = 'B02号直流'
$data[14,2] = 45935.026076388887
$data[14,3] = 45936.347141203703
$data[14,4] = 31.705555555585306
$data[15,0] = '长沙特来电飞狐四方坪南区充电站'
$data[15,1] = '101号直流'
$data[15,2] = 45935.0465625
$data[15,3] = 45936.347141203703
$data[15,4] = 31.213888888887595
$data[16,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[16,1] = '208号直流'
$data[16,2] = 45935.057013888887
$data[16,3] = 45936.347141203703
$data[16,4] = 30.963055555592291
$data[17,0] = '长沙特来电飞狐四方坪西区充电站'
$data[17,1] = '501号直流'
$data[17,2] = 45935.118333333332
$data[17,3] = 45936.347141203703
$data[17,4] = 29.491388888913207
$data[18,0] = '长沙特来电飞狐四方坪西区充电站'
$data[18,1] = '205号直流'
$data[18,2] = 45935.183657407404
$data[18,3] = 45936.347141203703
$data[18,4] = 27.923611111182254
$data[19,0] = '长沙特来电飞狐四方坪西区充电站'
$data[19,1] = '705号直流'
$data[19,2] = 45935.199178240742
$data[19,3] = 45936.347141203703
$data[19,4] = 27.551111111068167
$data[20,0] = '长沙特来电飞狐四方坪南区充电站'
$data[20,1] = '304号直流'
$data[20,2] = 45935.232499999998
$data[20,3] = 45936.347141203703
$data[20,4] = 26.75138888892252
$data[21,0] = '长沙特来电飞狐四方坪南区充电站'
$data[21,1] = '201号直流'
$data[21,2] = 45935.236712962964
$data[21,3] = 45936.347141203703
$data[21,4] = 26.650277777749579
$data[22,0] = '长沙特来电飞狐四方坪西区充电站'
$data[22,1] = '404号直流'
$data[22,2] = 45935.240381944444
$data[22,3] = 45936.347141203703
$data[22,4] = 26.562222222215496
$data[23,0] = '长沙特来电飞狐四方坪西区充电站'
$data[23,1] = '403号直流'
$data[23,2] = 45935.245081018518
$data[23,3] = 45936.347141203703
$data[23,4] = 26.449444444442634
$data[24,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[24,1] = '107号直流'
$data[24,2] = 45935.359930555554
$data[24,3] = 45936.347141203703
$data[24,4] = 23.693055555573665
$data[25,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[25,1] = '109号直流'
$data[25,2] = 45935.367847222224
$data[25,3] = 45936.347141203703
$data[25,4] = 23.503055555513129
$data[26,0] = '长沙特来电飞狐四方坪西区充电站'
$data[26,1] = '604号直流'
$data[26,2] = 45935.376030092593
$data[26,3] = 45936.347141203703
$data[26,4] = 23.306666666641831
$data[27,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[27,1] = '104号直流'
$data[27,2] = 45935.455196759256
$data[27,3] = 45936.347141203703
$data[27,4] = 21.406666666734964
$data[28,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[28,1] = '304号直流'
$data[28,2] = 45935.496365740742
$data[28,3] = 45936.347141203703
$data[28,4] = 20.418611111061182
$data[29,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[29,1] = '203号直流'
$data[29,2] = 45935.510324074072
$data[29,3] = 45936.347141203703
$data[29,4] = 20.083611111156642
$data[30,0] = '长沙特来电飞狐四方坪南区充电站'
$data[30,1] = '306号直流'
$data[30,2] = 45935.517800925925
$data[30,3] = 45936.347141203703
$data[30,4] = 19.904166666674428
$data[31,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[31,1] = '206号直流'
$data[31,2] = 45935.528865740744
$data[31,3] = 45936.347141203703
$data[31,4] = 19.638611111033242
$data[32,0] = '长沙特来电飞狐四方坪西区充电站'
$data[32,1] = '402号直流'
$data[32,2] = 45935.53224537037
$data[32,3] = 45936.347141203703
$data[32,4] = 19.557499999995343
$data[33,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[33,1] = '308号直流'
$data[33,2] = 45935.532754629632
$data[33,3] = 45936.347141203703
$data[33,4] = 19.545277777709998
$data[34,0] = '长沙特来电飞狐四方坪南区充电站'
$data[34,1] = '401号直流'
$data[34,2] = 45935.540162037039
$data[34,3] = 45936.347141203703
$data[34,4] = 19.367499999934807
$data[35,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[35,1] = '212号直流'
$data[35,2] = 45935.547060185185
$data[35,3] = 45936.347141203703
$data[35,4] = 19.201944444444962
$data[36,0] = '长沙特来电飞狐四方坪西区充电站'
$data[36,1] = '905号直流'
$data[36,2] = 45935.549293981479
$data[36,3] = 45936.347141203703
$data[36,4] = 19.148333333374467
$data[37,0] = '长沙特来电飞狐四方坪西区充电站'
$data[37,1] = '804号直流'
$data[37,2] = 45935.560706018521
$data[37,3] = 45936.347141203703
$data[37,4] = 18.874444444372784
$data[38,0] = '长沙特来电飞狐四方坪东区充电站'
$data[38,1] = '402号直流'
$data[38,2] = 45935.573391203703
$data[38,3] = 45936.347141203703
$data[38,4] = 18.570000000006985
$data[39,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[39,1] = '106号直流'
$data[39,2] = 45935.577557870369
$data[39,3] = 45936.347141203703
$data[39,4] = 18.470000000030268
$data[40,0] = '长沙特来电飞狐四方坪东区充电站'
$data[40,1] = '906号直流'
$data[40,2] = 45935.584687499999
$data[40,3] = 45936.347141203703
$data[40,4] = 18.29888888890855
$data[41,0] = '长沙特来电飞狐四方坪南区充电站'
$data[41,1] = '105号直流'
$data[41,2] = 45935.591550925928
$data[41,3] = 45936.347141203703
$data[41,4] = 18.134166666597594
$data[42,0] = '长沙特来电飞狐四方坪东区充电站'
$data[42,1] = '006B号直流'
$data[42,2] = 45935.591967592591
$data[42,3] = 45936.347141203703
$data[42,4] = 18.124166666704696
$data[43,0] = '长沙特来电飞狐四方坪东区充电站'
$data[43,1] = '905号直流'
$data[43,2] = 45935.597280092596
$data[43,3] = 45936.347141203703
$data[43,4] = 17.996666666585952
$data[44,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[44,1] = '312号直流'
$data[44,2] = 45935.627534722225
$data[44,3] = 45936.347141203703
$data[44,4] = 17.270555555471219
$data[45,0] = '长沙特来电飞狐四方坪西区充电站'
$data[45,1] = '903号直流'
$data[45,2] = 45935.648599537039
$data[45,3] = 45936.347141203703
$data[45,4] = 16.764999999955762
$data[46,0] = '长沙特来电飞狐四方坪西区充电站'
$data[46,1] = '805号直流'
$data[46,2] = 45935.662916666668
$data[46,3] = 45936.347141203703
$data[46,4] = 16.421388888848014
$data[47,0] = '长沙特来电飞狐四方坪西区充电站'
$data[47,1] = '904号直流'
$data[47,2] = 45935.670243055552
$data[47,3] = 45936.347141203703
$data[47,4] = 16.245555555622559
$data[48,0] = '长沙特来电飞狐四方坪西区充电站'
$data[48,1] = '405号直流'
$data[48,2] = 45935.678483796299
$data[48,3] = 45936.347141203703
$data[48,4] = 16.047777777712326
$data[49,0] = '长沙特来电飞狐四方坪南区充电站'
$data[49,1] = '206号直流'
$data[49,2] = 45935.708032407405
$data[49,3] = 45936.347141203703
$data[49,4] = 15.338611111161299
$data[50,0] = '长沙特来电飞狐四方坪西区充电站'
$data[50,1] = '702号直流'
$data[50,2] = 45935.710706018515
$data[50,3] = 45936.347141203703
$data[50,4] = 15.274444444512483
$data[51,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[51,1] = '103号直流'
$data[51,2] = 45935.721226851849
$data[51,3] = 45936.347141203703
$data[51,4] = 15.021944444510154
$data[52,0] = '长沙特来电飞狐四方坪东区充电站'
$data[52,1] = '003B号直流'
$data[52,2] = 45935.728483796294
$data[52,3] = 45936.347141203703
$data[52,4] = 14.8477777778171
$data[53,0] = '长沙特来电飞狐四方坪东区充电站'
$data[53,1] = '011A号直流'
$data[53,2] = 45935.745150462964
$data[53,3] = 45936.347141203703
$data[53,4] = 14.44777777773561
$data[54,0] = '长沙特来电飞狐四方坪西区充电站'
$data[54,1] = 'A05号直流'
$data[54,2] = 45935.793541666666
$data[54,3] = 45936.347141203703
$data[54,4] = 13.286388888896909
$data[55,0] = '长沙市开福区高岭香江国际城充电站建设项目'
$data[55,1] = '305号直流'
$data[55,2] = 45935.804583333331
$data[55,3] = 45936.347141203703
$data[55,4] = 13.021388888941146
$data[56,0] = '长沙特来电飞狐四方坪西区充电站'
$data[56,1] = 'A03号直流'
$data[56,2] = 45935.83017361111
$data[56,3] = 45936.347141203703
$data[56,4] = 12.407222222245764
$data[57,0] = '长沙特来电飞狐四方坪东区充电站'
$data[57,1] = '603号直流'
$data[57,2] = 45935.841516203705
$data[57,3] = 45936.347141203703
$data[57,4] = 12.134999999951106
$ws.Range("A2:E59").Value = $data

$ws.Range("H7").Select()

